# Word COM-interop script reproducing the commit:
#   "Have working test code including an applet that takes an input
#    string from a host and can return it later."
#
# Structural changes to word/document.xml:
#   1. The "_GoBack" bookmark moves from the last paragraph ("Get
#      OpenCard framework") to the (previously empty) paragraph right
#      after the Oracle hyperlink paragraph.
#   2. The text of the last paragraph changes from
#        "Get OpenCard framework"
#      to
#        "Is secure channel necessary in gpshell? What is the purpose."
#      while keeping the existing spell-check run split around the
#      single "proper noun" word (OpenCard -> gpshell) intact.

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -----------------------------
# Word only ever keeps a single "_GoBack" bookmark; adding a new one
# named "_GoBack" implicitly removes it from wherever it used to be
# (the last paragraph) and places it at the new range (the empty
# paragraph directly under the hyperlink).
$goBackPara = $d.Paragraphs.Item(3)
$d.Bookmarks.Add("_GoBack", $goBackPara.Range) | Out-Null

# --- 2. Rewrite the "Get OpenCard framework" paragraph ---------------
# Each replacement stays fully inside the run that already carries the
# text, so the existing run / proofErr (spell-check) split around the
# word that changes ("OpenCard" -> "gpshell") is preserved exactly.
$d.Content.Find.Execute("Get ", $true, $false, $false, $false, $false, $true, 1, $false, "Is secure channel necessary in ", 2) | Out-Null

$d.Content.Find.Execute("OpenCard", $true, $false, $false, $false, $false, $true, 1, $false, "gpshell", 2) | Out-Null

$d.Content.Find.Execute(" framework", $true, $false, $false, $false, $false, $true, 1, $false, "? What is the purpose.", 2) | Out-Null
